# #12 Boxes colors linked with text on issues slides
#
# Slide 7 ("Hard to deploy") has two thin rectangles that used to share the
# same scheme color (tx2, lumMod 60%/lumOff 40%). Recolor them to solid
# green / yellow and make the matching words in the body text bold and
# colored the same way, so the boxes are visually linked to the text.

function HexToRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$green  = HexToRGB "00B050"
$yellow = HexToRGB "FFFF00"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# --- Rectangle 5 (first thin box) -> solid green fill/outline -------------
$rectGreen = $s.Shapes.Item(1)
$rectGreen.Fill.ForeColor.RGB = $green
$rectGreen.Line.ForeColor.RGB = $green

# --- Rectangle 3 (second thin box) -> solid yellow fill/outline -----------
$rectYellow = $s.Shapes.Item(2)
$rectYellow.Fill.ForeColor.RGB = $yellow
$rectYellow.Line.ForeColor.RGB = $yellow

# --- Body text placeholder: bold + color the matching words ---------------
$content = $s.Shapes.Item(4)
$tr = $content.TextFrame.TextRange

# "... every support" -> "... every operating system" (bold, yellow)
$full = $tr.Text
$idx = $full.IndexOf(" support")
$word = $tr.Characters($idx + 2, 7)
$word.Text = "operating system"
$run = $tr.Characters($idx + 2, 16)
$run.Font.Bold = $true
$run.Font.Color.RGB = $yellow

# "... interpreting inputs" -> bold, green "inputs"
$full = $tr.Text
$idx = $full.IndexOf(" inputs")
$run2 = $tr.Characters($idx + 2, 6)
$run2.Font.Bold = $true
$run2.Font.Color.RGB = $green
